$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Change Log")
$ws.Activate()

# Row 7: new log entry dated 8/17/2025 (serial 45886), with a new
# "Changes" note (B7) and a new "Notes" note (D7). Column C7 (Editor "NW")
# is already correct and left untouched.
$ws.Range("A7").Value = 45886

$ws.Range("B7").Value = "Changes`n- MODIFIED: alu_control.vhd, alu.vhd`n- ADDED: AND, OR, XOR, Shift right and left logical, and shift right arithmetic                                                                                                                                                                                                                                        "

$ws.Range("D7").Value = "Notes`n- Haven't tested what you added to alu_control or alu yet, but it all compiles so far`nBugs`n- "

# Update the sheet's saved view state: scrolled so row 4 is at the top,
# with D8 as the active/selected cell.
$ws.Range("D8").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
